# Saldo.xlsx update
#
# The "Export" sheet is a flat dump of Conta/Nome/Saldo rows (row 1 is the
# header). This change:
#   1. Removes the RODOLFO row (004213929, 120000) entirely.
#   2. Removes the JOSE row (002687737, 12000) entirely.
#   3. Moves the ANGELICA row (004207374) down so it lands right after the
#      DANIELI row (004377713), and corrects her balance from 126450.48
#      down to 450.48.
# Everything else keeps its original relative order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1 & 2: delete the RODOLFO and JOSE rows -------------------------------
# Row 5 = 002687737 / JOSE / 12000
# Row 3 = 004213929 / RODOLFO / 120000
# Delete the lower row first so the higher row's index doesn't shift.
$ws.Rows(5).Delete()
$ws.Rows(3).Delete()

# --- 3: relocate ANGELICA ----------------------------------------------------
# After the two deletions above, ANGELICA is still row 2 (004207374 / ANGELICA
# / 126450.48), and DANIELI (004377713) has shifted from row 22 up to row 20.
$angelicaAccount = $ws.Cells.Item(2, 1).Value2
$angelicaName = $ws.Cells.Item(2, 2).Value2

$ws.Rows(2).Delete()

# DANIELI is now at row 19 (20 - 1, from the ANGELICA delete above); insert
# the new row right after it, at row 20.
$ws.Rows(20).Insert()
$ws.Cells.Item(20, 1).NumberFormat = "@"
$ws.Cells.Item(20, 1).Value = $angelicaAccount
$ws.Cells.Item(20, 2).Value = $angelicaName
$ws.Cells.Item(20, 3).Value = 450.48
